$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.434.66"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "1.463.25"
$ws.Range("E3").Value = "  +3.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  -5.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3649"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3072"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.034"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06558"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9968"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.395"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.115"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001023"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "1.458.57"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9586"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05748"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.418"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.237"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").Value = "20.435.56"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.082"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "1.610.26"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.816"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.865"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.87%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7872"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07788"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.504"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05700"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.663"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02026"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9491"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1859"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.393"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5256"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.743"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06422"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9882"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
